$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values for rows 2-5 (new dataset of 1000 records / updated sample)
# Row 2
$ws.Range("A2").Value = 45167.50694444445
$ws.Range("B2").Value = 19.217
$ws.Range("C2").Value = 12.901
$ws.Range("D2").Value = 4.042
$ws.Range("E2").Value = 40.812
$ws.Range("F2").Value = 32.818
$ws.Range("G2").Value = 15.123
$ws.Range("H2").Value = 47.986
$ws.Range("I2").Value = 23.269
$ws.Range("J2").Value = 9.710000000000001
$ws.Range("K2").Value = 14.67
$ws.Range("L2").Value = 16.076
$ws.Range("M2").Value = 16.742
$ws.Range("N2").Value = 4.827
$ws.Range("O2").Value = 15.038
$ws.Range("P2").Value = 20.994
$ws.Range("Q2").Value = 12.85
$ws.Range("R2").Value = 3.46
$ws.Range("S2").Value = 2.249
$ws.Range("T2").Value = 221.547
$ws.Range("U2").Value = 41.81
$ws.Range("V2").Value = 13.881
$ws.Range("W2").Value = 27.553
$ws.Range("X2").Value = 14.055
$ws.Range("Y2").Value = 3.03
$ws.Range("Z2").Value = 24.312
$ws.Range("AA2").Value = 12.261
$ws.Range("AB2").Value = 11.125
$ws.Range("AC2").Value = 13.047
$ws.Range("AD2").Value = 16.565
$ws.Range("AE2").Value = 3.456
$ws.Range("AF2").Value = 42.557
$ws.Range("AG2").Value = 7.647
$ws.Range("AH2").Value = 17.354

# Row 3
$ws.Range("A3").Value = 45167.51388888889
$ws.Range("B3").Value = 7.687
$ws.Range("C3").Value = 5.061
$ws.Range("D3").Value = 1.518
$ws.Range("E3").Value = 16.474
$ws.Range("F3").Value = 13.027
$ws.Range("G3").Value = 6.05
$ws.Range("H3").Value = 26.279
$ws.Range("I3").Value = 9.308
$ws.Range("J3").Value = 3.867
$ws.Range("K3").Value = 5.615
$ws.Range("L3").Value = 6.57
$ws.Range("M3").Value = 6.766
$ws.Range("N3").Value = 1.937
$ws.Range("O3").Value = 6.015
$ws.Range("P3").Value = 8.375999999999999
$ws.Range("Q3").Value = 5.394
$ws.Range("R3").Value = 1.481
$ws.Range("S3").Value = 0.8070000000000001
$ws.Range("T3").Value = 84.26000000000001
$ws.Range("U3").Value = 16.997
$ws.Range("V3").Value = 5.552
$ws.Range("W3").Value = 11.042
$ws.Range("X3").Value = 5.631
$ws.Range("Y3").Value = 1.366
$ws.Range("Z3").Value = 12.426
$ws.Range("AA3").Value = 4.904
$ws.Range("AB3").Value = 4.582
$ws.Range("AC3").Value = 5.349
$ws.Range("AD3").Value = 6.71
$ws.Range("AE3").Value = 1.265
$ws.Range("AF3").Value = 24.163
$ws.Range("AG3").Value = 2.975
$ws.Range("AH3").Value = 6.943

# Row 4
$ws.Range("A4").Value = 45167.52083333334
$ws.Range("B4").Value = 16.815
$ws.Range("C4").Value = 12.169
$ws.Range("D4").Value = 1.335
$ws.Range("E4").Value = 36.463
$ws.Range("F4").Value = 29.696
$ws.Range("G4").Value = 13.233
$ws.Range("H4").Value = 49.061
$ws.Range("I4").Value = 20.36
$ws.Range("J4").Value = 8.919
$ws.Range("K4").Value = 13.187
$ws.Range("L4").Value = 14.63
$ws.Range("M4").Value = 15.348
$ws.Range("N4").Value = 4.226
$ws.Range("O4").Value = 13.159
$ws.Range("P4").Value = 18.62
$ws.Range("Q4").Value = 11.225
$ws.Range("R4").Value = 1.054
$ws.Range("S4").Value = 0.83
$ws.Range("T4").Value = 192.956
$ws.Range("U4").Value = 36.731
$ws.Range("V4").Value = 12.146
$ws.Range("W4").Value = 24.53
$ws.Range("X4").Value = 12.841
$ws.Range("Y4").Value = 2.125
$ws.Range("Z4").Value = 24.029
$ws.Range("AA4").Value = 10.728
$ws.Range("AB4").Value = 9.609
$ws.Range("AC4").Value = 11.279
$ws.Range("AD4").Value = 15.259
$ws.Range("AE4").Value = 0.773
$ws.Range("AF4").Value = 44.353
$ws.Range("AG4").Value = 6.779
$ws.Range("AH4").Value = 15.185

# Row 5
$ws.Range("A5").Value = 45167.52777777778
$ws.Range("B5").Value = 6.73
$ws.Range("C5").Value = 4.71
$ws.Range("D5").Value = 0.78
$ws.Range("E5").Value = 14.56
$ws.Range("F5").Value = 11.68
$ws.Range("G5").Value = 5.29
$ws.Range("H5").Value = 24.21
$ws.Range("I5").Value = 8.140000000000001
$ws.Range("J5").Value = 3.52
$ws.Range("K5").Value = 5.08
$ws.Range("L5").Value = 5.85
$ws.Range("M5").Value = 6.09
$ws.Range("N5").Value = 1.7
$ws.Range("O5").Value = 5.26
$ws.Range("P5").Value = 7.43
$ws.Range("Q5").Value = 4.62
$ws.Range("R5").Value = 0.73
$ws.Range("S5").Value = 0.43
$ws.Range("T5").Value = 72.81999999999999
$ws.Range("U5").Value = 14.9
$ws.Range("V5").Value = 4.86
$ws.Range("W5").Value = 9.84
$ws.Range("X5").Value = 5.08
$ws.Range("Y5").Value = 1
$ws.Range("Z5").Value = 11.33
$ws.Range("AA5").Value = 4.29
$ws.Range("AB5").Value = 3.92
$ws.Range("AC5").Value = 4.59
$ws.Range("AD5").Value = 6.07
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 22.25
$ws.Range("AG5").Value = 2.65
$ws.Range("AH5").Value = 6.08

# Remove row 6 (dataset now has one fewer displayed row)
$ws.Rows.Item(6).Delete()

# Adjust column widths: change width 7 -> 8 for specific columns
# Offset of 5/6 (0.8333...) accounts for COM ColumnWidth vs raw XML width padding
$ws.Range("C1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("G1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("K1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("L1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("O1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("Q1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("V1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("X1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AA1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AB1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AC1").EntireColumn.ColumnWidth = 7.166666666666667
